# Week 15 simulations update
# Appends new simulated-game data points to the running per-play stat
# lists on the YDS and ST sheets, and refreshes the derived totals on the
# OFF / DEF / ST / TURNS / PEN summary sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append new week's numbers to the space-separated run lists
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = $wsYDS.Range("B2").Value2 + " " + "5 -1 0 2 6 5 4 11 0 4 -1 4 11 4 0 11 1 4 3 11 8 5 9 0 2 1 -2 1 11 11 4 2 3 5 3"
$wsYDS.Range("B3").Value = $wsYDS.Range("B3").Value2 + " " + "-2 6 36 9 13 3 20 4 4 14 12 7 14 3 11 6 2 23"
$wsYDS.Range("C2").Value = $wsYDS.Range("C2").Value2 + " " + "3 -1 11 8 -1 12 8 4 2 4 1 4 1 4 -5 -1 5 11 17 4 0 1 2 -2 2 4"
$wsYDS.Range("C3").Value = $wsYDS.Range("C3").Value2 + " " + "7 11 15 22 10 23 9 8 14 6 18 7 9 3 19 19 13 18 11"

# ---------------------------------------------------------------------
# OFF sheet: updated season totals
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("B2").Value = 9
$wsOFF.Range("C2").Value = 368
$wsOFF.Range("E2").Value = 17
$wsOFF.Range("F2").Value = 87
$wsOFF.Range("G2").Value = 83
$wsOFF.Range("J2").Value = 47
$wsOFF.Range("O2").Value = 42
$wsOFF.Range("P2").Value = 22

$wsOFF.Range("B3").Value = 12
$wsOFF.Range("C3").Value = 294
$wsOFF.Range("D3").Value = 8
$wsOFF.Range("E3").Value = 56
$wsOFF.Range("F3").Value = 226
$wsOFF.Range("H3").Value = 56
$wsOFF.Range("I3").Value = 123
$wsOFF.Range("L3").Value = 517
$wsOFF.Range("M3").Value = 350
$wsOFF.Range("Q3").Value = 923

# ---------------------------------------------------------------------
# DEF sheet: updated season totals
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 336
$wsDEF.Range("E2").Value = 16
$wsDEF.Range("F2").Value = 86
$wsDEF.Range("G2").Value = 111
$wsDEF.Range("I2").Value = 10
$wsDEF.Range("J2").Value = 47
$wsDEF.Range("N2").Value = 20
$wsDEF.Range("O2").Value = 39
$wsDEF.Range("P2").Value = 17

$wsDEF.Range("B3").Value = 18
$wsDEF.Range("C3").Value = 357
$wsDEF.Range("E3").Value = 64
$wsDEF.Range("F3").Value = 202
$wsDEF.Range("I3").Value = 101
$wsDEF.Range("J3").Value = 103
$wsDEF.Range("L3").Value = 564
$wsDEF.Range("M3").Value = 384
$wsDEF.Range("Q3").Value = 959

# ---------------------------------------------------------------------
# ST sheet: updated season totals + appended kicker/return distance lists
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 141
$wsST.Range("D2").Value = 110
$wsST.Range("F2").Value = 85
$wsST.Range("G2").Value = 81
$wsST.Range("H2").Value = 6
$wsST.Range("J2").Value = 56
$wsST.Range("K2").Value = 54
$wsST.Range("L2").Value = 26
$wsST.Range("M2").Value = 19

$wsST.Range("B3").Value = 74

$wsST.Range("B4").Value = $wsST.Range("B4").Value2 + " " + "63 66 64 66"
$wsST.Range("B5").Value = $wsST.Range("B5").Value2 + " " + "26 28 0 25"
$wsST.Range("B6").Value = $wsST.Range("B6").Value2 + " " + "38"
$wsST.Range("D3").Value = $wsST.Range("D3").Value2 + " " + "46 64"
$wsST.Range("D4").Value = $wsST.Range("D4").Value2 + " " + "0 0"
$wsST.Range("D5").Value = "0 0 0 0"

# ---------------------------------------------------------------------
# TURNS sheet: updated Road totals
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("C3").Value = 13
$wsTURNS.Range("D3").Value = 19
$wsTURNS.Range("E3").Value = 16

# ---------------------------------------------------------------------
# PEN sheet: updated OFF penalty totals
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 36
$wsPEN.Range("B3").Value = 20
